# Insert a new weekly price record as row 63 (Acelga, "Segunda" quality,
# reported 2023-01-17 / serial 44943), pushing the existing rows 63-84
# down to 64-85. The new row reuses the same market/category/origin
# metadata as the row that follows it, differing only in date, quality
# and volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63:84 down to 64:85, leaving a blank row 63.
$ws.Rows.Item(63).Insert()

# Seed the new row with the (now-shifted) following row's data, then
# overwrite the fields that differ for this record.
$ws.Range("A64:R64").Copy($ws.Range("A63:R63"))

$ws.Cells.Item(63, 4).Value = 44943      # Fecha
$ws.Cells.Item(63, 9).Value = "Segunda"  # Calidad
$ws.Cells.Item(63, 10).Value = 300       # Volumen

Write-Output "Inserted new row 63 for Acelga weekly update"
